# Update "想去人数" (want-to-go count) values in column F on the sheets
# that hold the full conference listing: "展览" and "全部类型".
# Row -> (old value, new value):
#   3  : 2195  -> 2196
#   5  : 13049 -> 13051
#   10 : 1175  -> 1176
#   11 : 975   -> 976
#   25 : 5375  -> 5378
#   29 : 13    -> 14

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 2196
    5  = 13051
    10 = 1176
    11 = 976
    25 = 5378
    29 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
